# Renaming config and code variables for better understanding
#
# - Renames several sheet tabs to clearer names.
# - Renames a handful of header / label cells to match the new
#   "reward" / "delay" terminology (replacing "prize" / "lane-time").
# - Moves the active-sheet/tab-selection + cell-selection state around
#   (selection now lives on the first sheet instead of Debug).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Rename worksheets (by their CURRENT names, before any renaming
#    changes what "Item(...)" would resolve to).
# ---------------------------------------------------------------------
$wsTrainings   = $wb.Worksheets.Item("Trainings")
$wsExperiments = $wb.Worksheets.Item("Experiments")
$wsTexts       = $wb.Worksheets.Item("Texts")
$wsPrizes      = $wb.Worksheets.Item("Prizes")
$wsTimes       = $wb.Worksheets.Item("Times")
$wsDebug       = $wb.Worksheets.Item("Debug")
$wsExport      = $wb.Worksheets.Item("Export")

$wsTrainings.Name   = "Practices"
$wsExperiments.Name = "Tasks"
$wsPrizes.Name      = "Rewards"
$wsTimes.Name       = "Delays"
# Texts, Debug, Export keep their names.

# ---------------------------------------------------------------------
# 2. Rename header / label cell text.
# ---------------------------------------------------------------------

# "Practices" (was Trainings) header row B1:D1
$wsTrainings.Range("B1").Value = "Immediate Reward Value"
$wsTrainings.Range("C1").Value = "Delayed Reward Lane"
$wsTrainings.Range("D1").Value = "Delayed Reward Value"

# "Tasks" (was Experiments) header row B1:D1
$wsExperiments.Range("B1").Value = "Immediate Reward Value"
$wsExperiments.Range("C1").Value = "Delayed Reward Lane"
$wsExperiments.Range("D1").Value = "Delayed Reward Value"

# "Delays" (was Times) header row A1:B1
$wsTimes.Range("A1").Value = "Lane"
$wsTimes.Range("B1").Value = "Time"

# ---------------------------------------------------------------------
# 3. Update sheet view / selection state.
#    Selection ends up active on the "Practices" sheet (tabSelected),
#    matching the diff; Debug loses its tabSelected flag.
# ---------------------------------------------------------------------

$wsDebug.Activate()
$wsDebug.Range("C3").Select()

$wsExperiments.Activate()
$wsExperiments.Range("G6").Select()

$wsTimes.Activate()
$wsTimes.Range("G9").Select()

$wsTrainings.Activate()
$wsTrainings.Range("E16").Select()
